$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Row 9: "Update thunder data analysis..." task is now done.
# Clear the Comments cell (G9) and set the Done date (F9) instead.
$ws.Range("G9").ClearContents()
$ws.Range("F8").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Value = 45936

# Row 10: "Convert Meteorite analysis scripts from Jupyter to .py" assigned to Sprint 4.
$ws.Range("E10").Value = 4

# New row 16: add a new backlog task.
$ws.Range("A16").Value = "Handle Leap years in Thunder Average Script"
$ws.Range("B16").Value = "$"
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = "No"

# Copy formatting from an existing similar row (row 15) onto the new row's used cells.
$ws.Range("A15:D15").Copy() | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null

# Update the selected cell on the sheet.
$ws.Range("G8").Select() | Out-Null
